$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.355.91"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +4.16%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.347.53"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +2.64%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D5').Value = "'547.13"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +3.18%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'132.50"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.41%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'2.343.80"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +2.58%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  +2.21%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  +1.25%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +1.02%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.334"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +1.72%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'23.91"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +2.10%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'2.763.74"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.62%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'60.306.57"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +4.20%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.0000133"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +1.67%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'2.344.92"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +1.69%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'10.65"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +1.43%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.41%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  +7.12%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'314.18"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +1.21%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.12%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'63.55"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +1.90%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'0.173"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +3.25%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.07%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.88%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +8.33%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +2.83%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'171.59"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.62%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  +12.88%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +2.20%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +4.11%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +14.50%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  +1.11%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'18.08"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +2.12%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.03%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.05%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +7.69%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'325.87"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +13.86%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'38.13"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -0.94%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  +3.18%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'141.49"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.31%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  +1.76%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.0948"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +0.26%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'19.57"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +9.28%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.05%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +1.99%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  +2.31%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'11.04"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +1.17%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  +13.70%  "
$ws.Range('E51').Style = 'Normal'
